$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

$ws.Range("A2").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B2").Value = "Prejuveniles"
$ws.Range("E2").Value = "Garrone, Federico Daniel"
$ws.Range("F2").Value = 90
